$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize dye-name text in column A (uppercase, cleanup suffixes/typos)
$ws.Range("A2").Value = '19TH!'
$ws.Range("A3").Value = '6 OF EVERYTHING'
$ws.Range("A5").Value = 'ABSOLUTE MAGENTA'
$ws.Range("A6").Value = 'AFTER HOURS'
$ws.Range("A7").Value = 'ALIVE'
$ws.Range("A8").Value = 'ALL AMERICAN HEN'
$ws.Range("A9").Value = 'AMBER GLASS'
$ws.Range("A11").Value = 'AMETHYST INK'
$ws.Range("A13").Value = 'ANGELICA'
$ws.Range("A14").Value = 'ANGEL’S LANDING'
$ws.Range("A15").Value = 'ANYTHING GOES'
$ws.Range("A16").Value = 'APATHY'
$ws.Range("A17").Value = 'AQUA JET'
$ws.Range("A18").Value = 'AS A BIRD'
$ws.Range("A19").Value = 'ATOMIC BLUE'
$ws.Range("A20").Value = 'AWESOME! WOW!'
$ws.Range("A21").Value = 'AZURE COVE'
$ws.Range("A23").Value = 'BEDTIME'
$ws.Range("A24").Value = 'BELLE EPOQUE'
$ws.Range("A25").Value = 'BELOW HORIZON'
$ws.Range("A26").Value = 'BERMUDA TEAL'
$ws.Range("A28").Value = 'BITTER MALAISE'
$ws.Range("A30").Value = 'BLACK PEARL'
$ws.Range("A31").Value = 'BLANCHEFLEUR'
$ws.Range("A33").Value = 'BLUE LAGOON'
$ws.Range("A34").Value = 'BLUE SAGE'
$ws.Range("A35").Value = 'BLUE SULK'
$ws.Range("A36").Value = 'BOOT CAMP'
$ws.Range("A40").Value = 'CABARET'
$ws.Range("A41").Value = 'CALDERA'
$ws.Range("A42").Value = 'CALLOUS PINK'
$ws.Range("A43").Value = 'CANDY LOVE'
$ws.Range("A44").Value = 'CARMEN MIRANDA'
$ws.Range("A46").Value = 'CAROUSEL'
$ws.Range("A47").Value = 'CATCHING STARS'
$ws.Range("A48").Value = 'CEDAR CREEK'
$ws.Range("A49").Value = 'CHARCOAL PRISMATIC'
$ws.Range("A50").Value = 'CHARGED CHERRY'
$ws.Range("A51").Value = 'CHARLI AU LAIT'
$ws.Range("A52").Value = 'CHESLEY'
$ws.Range("A54").Value = 'CHOCOLATE NIGHT'
$ws.Range("A55").Value = 'CINNAMON GIRL'
$ws.Range("A56").Value = 'CITIES TOWERS & BRIDGES'
$ws.Range("A57").Value = 'CLOUD TO GROUND'
$ws.Range("A63").Value = 'CRYING DOVE'
$ws.Range("A65").Value = 'DAMP PILLOW'
$ws.Range("A68").Value = 'DEEP REGRET'
$ws.Range("A69").Value = 'DEEP SEAFLOWER'
$ws.Range("A70").Value = 'DESERT CITY WATTAGE'
$ws.Range("A71").Value = 'DESERT OASIS'
$ws.Range("A72").Value = 'DEVIL’S GARDEN'
$ws.Range("A74").Value = 'DIVINE'
$ws.Range("A75").Value = 'DO RE ME'
$ws.Range("A76").Value = 'DREAM ON!'
$ws.Range("A77").Value = 'DREAMGIRLS'
$ws.Range("A78").Value = 'DUNE'
$ws.Range("A79").Value = 'DUSK TO DAWN'
$ws.Range("A80").Value = 'EATEN THE PLUMS'
$ws.Range("A83").Value = 'ELEMENT 79'
$ws.Range("A84").Value = 'ELIZA'
$ws.Range("A86").Value = 'EMERALD DARKNESS'
$ws.Range("A88").Value = 'FABLE'
$ws.Range("A89").Value = 'FANCY'
$ws.Range("A90").Value = 'FATAL DUEL'
$ws.Range("A91").Value = 'FEDERAL BANK'
$ws.Range("A92").Value = 'FERN SLIPPER'
$ws.Range("A93").Value = 'FIERCE SCARLET'
$ws.Range("A94").Value = 'FLOWER DRUM SONG'
$ws.Range("A95").Value = 'FORGET ME'
$ws.Range("A97").Value = 'FORTUNE TELLER'
$ws.Range("A99").Value = 'GET DOWN GREENY'
$ws.Range("A100").Value = 'GG LAUREL'
$ws.Range("A103").Value = 'GOBLIN VALLEY'
$ws.Range("A104").Value = 'GOLD EXPERIENCE'
$ws.Range("A105").Value = 'GOLD MINE'
$ws.Range("A106").Value = 'GOOD LUCK JADE'
$ws.Range("A107").Value = 'GREEN LANTERN'
$ws.Range("A108").Value = 'GREY TABBY'
$ws.Range("A110").Value = 'GUTHRIE PEAK'
$ws.Range("A111").Value = 'HAYWIRE'
$ws.Range("A112").Value = 'HEAR IT HAWAII'
$ws.Range("A114").Value = 'HELLO DOLLY!'
$ws.Range("A115").Value = 'HER BEACON HAND'
$ws.Range("A116").Value = 'HER CHARISMA'
$ws.Range("A118").Value = 'HEY SUNSHINE!'
$ws.Range("A120").Value = 'HOT AIR BALLOON'
$ws.Range("A121").Value = 'HOT PANTS'
$ws.Range("A122").Value = 'HUMDRUM'
$ws.Range("A123").Value = 'ICELAND'
$ws.Range("A124").Value = 'ICY RECEPTION'
$ws.Range("A127").Value = 'IT’S A SPARKLER'
$ws.Range("A128").Value = 'IVY SNOWBELL'
$ws.Range("A129").Value = 'JAY FEATHER'
$ws.Range("A131").Value = 'JESSAMYN'
$ws.Range("A132").Value = 'JOCELYN'
$ws.Range("A133").Value = 'JOSHUA TREE'
$ws.Range("A134").Value = 'JULIETTE’S BLUSH'
$ws.Range("A135").Value = 'JUST TO SAY'
$ws.Range("A136").Value = 'K. CHAMELEON'
$ws.Range("A137").Value = 'KISS ME KATE'
$ws.Range("A138").Value = 'KYOTO SUNSET'
$ws.Range("A139").Value = 'LA LUZ'
$ws.Range("A140").Value = 'LAFAYETTE'
$ws.Range("A141").Value = 'LAVENDER BLOOM'
$ws.Range("A142").Value = 'LAY A ROSE'
$ws.Range("A143").Value = 'LEAF SEED BEAN'
$ws.Range("A144").Value = 'LEPRECHAUN BALLET'
$ws.Range("A145").Value = 'LIBERATION'
$ws.Range("A146").Value = 'LIFE OF THE GODS'
$ws.Range("A147").Value = 'LIPSTICK LAVA'
$ws.Range("A148").Value = 'LITTLE GREEN MEN'
$ws.Range("A149").Value = 'LOOSE GEMS'
$ws.Range("A150").Value = 'LOST IN PLUM'
$ws.Range("A151").Value = 'MAGIC ORCHID'
$ws.Range("A152").Value = 'MALIBU SAIL'
$ws.Range("A153").Value = 'MEDIEVAL'
$ws.Range("A154").Value = 'MELON BOMB'
$ws.Range("A155").Value = 'MERCADO LIGHTS'
$ws.Range("A156").Value = 'MERMAID SHOES'
$ws.Range("A157").Value = 'MESA'
$ws.Range("A158").Value = 'MIAMI RED'
$ws.Range("A159").Value = 'MILD TEDIUM'
$ws.Range("A161").Value = 'MIST OF NYX'
$ws.Range("A162").Value = 'MOCCASIN'
$ws.Range("A163").Value = 'MOD SQUAD'
$ws.Range("A165").Value = 'MY FAIR LADY'
$ws.Range("A166").Value = 'NAKED SHAME'
$ws.Range("A167").Value = 'NAVY ZEAL'
$ws.Range("A169").Value = 'NIGHTHAWKS'
$ws.Range("A170").Value = 'NO LIMIT'
$ws.Range("A171").Value = 'NOVEMBER MUSE'
$ws.Range("A172").Value = 'OLD PUEBLO'
$ws.Range("A173").Value = 'OLD VINE'
$ws.Range("A174").Value = 'OUT BEYOND'
$ws.Range("A176").Value = 'PEACEFUL TRANSITION'
$ws.Range("A177").Value = 'PEACOCK SHADOW'
$ws.Range("A178").Value = 'PENCHANT'
$ws.Range("A179").Value = 'PETAL SHOWER'
$ws.Range("A180").Value = 'PETRIFIED FOREST'
$ws.Range("A181").Value = 'PICKLE BALL '
$ws.Range("A182").Value = 'PINK BEAN'
$ws.Range("A183").Value = 'PINKY'
$ws.Range("A186").Value = 'POMAGRENADE'
$ws.Range("A188").Value = 'POWER PLANT'
$ws.Range("A189").Value = 'PRICKLY PEAR'
$ws.Range("A190").Value = 'PRINCE WILLIAM'
$ws.Range("A191").Value = 'PUNKY FUSCHIA'
$ws.Range("A192").Value = 'PURPLE RAIN'
$ws.Range("A193").Value = 'QUEEN’S LAKE'
$ws.Range("A194").Value = 'RABBIT EARS'
$ws.Range("A196").Value = 'RAIN ON ME'
$ws.Range("A197").Value = 'RASPBERRY BLAZE'
$ws.Range("A198").Value = 'RED RUSH'
$ws.Range("A199").Value = 'REGENCY'
$ws.Range("A200").Value = 'RELISH THE VOTE!'
$ws.Range("A202").Value = 'REYNOLDS'
$ws.Range("A203").Value = 'RING THE BLOSSOM BELL'
$ws.Range("A204").Value = 'RIO VERDE'
$ws.Range("A205").Value = 'RIOT GIRLS'
$ws.Range("A206").Value = 'ROCHAMBEAU'
$ws.Range("A207").Value = 'ROMEO BLUE'
$ws.Range("A208").Value = 'ROSALITA'
$ws.Range("A209").Value = 'ROSE ANGUISH'
$ws.Range("A211").Value = 'ROUTE 66'
$ws.Range("A213").Value = 'SCORCHED LIME'
$ws.Range("A214").Value = 'SECRET GARDEN'
$ws.Range("A215").Value = 'SHADOW BOX'
$ws.Range("A216").Value = 'SHE WALKS IN BEAUTY'
$ws.Range("A217").Value = 'SHINY MOSS'
$ws.Range("A218").Value = 'SHUYLER LAKE'
$ws.Range("A220").Value = 'SIMONE & SUSAN'
$ws.Range("A221").Value = 'SINGING IN THE RAIN'
$ws.Range("A222").Value = 'SKIPPING STONE'
$ws.Range("A223").Value = 'SKYROCKET'
$ws.Range("A224").Value = 'SLAMMIN LEMON'
$ws.Range("A226").Value = 'SOLAR GLITTER'
$ws.Range("A227").Value = 'SONORAN SKY'
$ws.Range("A228").Value = 'SOUTH PACIFIC'
$ws.Range("A229").Value = 'SPRING TICKLE'
$ws.Range("A230").Value = 'SPROUT'
$ws.Range("A231").Value = 'STARLESS SKY'
$ws.Range("A232").Value = 'STILL I RISE'
$ws.Range("A233").Value = 'STILL SPRUCE'
$ws.Range("A234").Value = 'STRANGE HARVEST'
$ws.Range("A235").Value = 'STRAW INTO GOLD'
$ws.Range("A236").Value = 'STRIKES TWICE'
$ws.Range("A237").Value = 'STRING & KEY'
$ws.Range("A238").Value = 'SUNDOWN ORCHID'
$ws.Range("A243").Value = 'TEX MEX'
$ws.Range("A244").Value = 'THE KING'
$ws.Range("A247").Value = 'TIFFANY BOX'
$ws.Range("A248").Value = 'TIME AWAY'
$ws.Range("A249").Value = 'TOKYO CREAM'
$ws.Range("A250").Value = 'TORCHWOOD'
$ws.Range("A252").Value = 'TRUE TO LIFE'
$ws.Range("A253").Value = 'TUCSON'
$ws.Range("A254").Value = 'TUMBLEWEED'
$ws.Range("A255").Value = 'TUSSIE MUSSIE'
$ws.Range("A256").Value = 'UNIMAGINABLE'
$ws.Range("A257").Value = 'UNIVERSE IS YOURS'
$ws.Range("A258").Value = 'UPTOWN ELECTRIC'
$ws.Range("A259").Value = 'URBAN FOSSIL'
$ws.Range("A260").Value = 'VAGUE UNEASE'
$ws.Range("A261").Value = 'VALENTINE'
$ws.Range("A264").Value = 'VIOLET’S BLUEBERRY'
$ws.Range("A265").Value = 'VISUAL PURPLE'
$ws.Range("A267").Value = 'WATCH MY RISING'
$ws.Range("A268").Value = 'WAY COOL CRISTAL'
$ws.Range("A269").Value = 'WHETSTONE'
$ws.Range("A270").Value = 'WHIP & CHILL'
$ws.Range("A271").Value = 'WHISPER'
$ws.Range("A272").Value = 'WICKED ROYAL'

# Remove frozen panes (B2 freeze no longer present in target)
$excel.ActiveWindow.FreezePanes = $false

# Row 117 height adjustment
$ws.Rows.Item(117).RowHeight = 25.85
